$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72") entirely; all rows below shift up by one,
# and the used range shrinks from A1:F63 to A1:F62.
$ws.Rows.Item(2).Delete()
